$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 48.204556939412974
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 51.680545304071345
$ws.Range("E2").Value = 56.084878227406591

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 46.533560905799668
$ws.Range("D3").Value = 45.332565871072198
$ws.Range("E3").Value = 54.100473278624925

# Update selection to match new selected range
$null = $ws.Range("B1:E3").Select()
